# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Rule R40 (row 11) is renumbered: its "Rule" label cell (B11) is changed
# from the text "R40" to the text "1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Enter "1" as literal text (not the number 1) in B11, exactly like typing
# '1 into the cell in Excel - the leading apostrophe forces a text entry
# even though it looks numeric.
$ws.Range("B11").Value = "'1"
